$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (class_id 0 - dress)
$ws.Range("C2").Value = 63
$ws.Range("D2").Value = 0.2758620689655172
$ws.Range("E2").Value = 0.253968253968254
$ws.Range("F2").Value = 0.2644628099173554

# Row 3 (class_id 1 - high_heel)
$ws.Range("C3").Value = 25
$ws.Range("D3").Value = 0.125
$ws.Range("E3").Value = 0.12
$ws.Range("F3").Value = 0.1224489795918367

# Row 4 (class_id 2 - handbag)
$ws.Range("C4").Value = 26
$ws.Range("D4").Value = 0.25
$ws.Range("E4").Value = 0.03846153846153846
$ws.Range("F4").Value = 0.06666666666666667

# Row 5 (class_id 3 - skirt)
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 0.1707317073170732
$ws.Range("E5").Value = 0.4666666666666667
$ws.Range("F5").Value = 0.25

# Row 6 (class_id 4 - outerwear)
$ws.Range("C6").Value = 56
$ws.Range("D6").Value = 0.1842105263157895
$ws.Range("E6").Value = 0.125
$ws.Range("F6").Value = 0.148936170212766

# Row 7 (class_id 5 - boot)
$ws.Range("C7").Value = 9
